# Auto-generated edit script: updates FFXIV crafting-leve market price/profit
# figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 740.5
$ws.Range("I96").Value = 276
$ws.Range("J96").Value = 1437.25
$ws.Range("K96").Value = 828
$ws.Range("L96").Value = 4311.75
$ws.Range("M96").Value = 545
$ws.Range("N96").Value = -7057.75
$ws.Range("H103").Value = 12445.429
$ws.Range("I103").Value = 353
$ws.Range("J103").Value = 85000
$ws.Range("K103").Value = 1059
$ws.Range("L103").Value = 255000
$ws.Range("M103").Value = -473
$ws.Range("N103").Value = -256172
$ws.Range("H106").Value = 2215.7144
$ws.Range("I106").Value = 1201.3334
$ws.Range("J106").Value = 2976.5
$ws.Range("K106").Value = 1201.3334
$ws.Range("L106").Value = 2976.5
$ws.Range("M106").Value = -570.3334
$ws.Range("N106").Value = -4238.5
$ws.Range("H129").Value = 87213
$ws.Range("J129").Value = 119773
$ws.Range("L129").Value = 359319
$ws.Range("N129").Value = -369319
$ws.Range("H138").Value = 2231.707
$ws.Range("I138").Value = 904.70966
$ws.Range("J138").Value = 2836.6619
$ws.Range("K138").Value = 2714.12898
$ws.Range("L138").Value = 8509.985700000001
$ws.Range("M138").Value = 2425.87102
$ws.Range("N138").Value = -18789.9857
$ws.Range("H141").Value = 12080.818
$ws.Range("I141").Value = 15048.625
$ws.Range("J141").Value = 4166.6665
$ws.Range("K141").Value = 45145.875
$ws.Range("L141").Value = 12499.9995
$ws.Range("M141").Value = -39965.875
$ws.Range("N141").Value = -22859.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6481
$ws.Range("I32").Value = 3788.1455
$ws.Range("J32").Value = 11966.444
$ws.Range("K32").Value = 3788.1455
$ws.Range("L32").Value = 11966.444
$ws.Range("M32").Value = -3501.1455
$ws.Range("N32").Value = -12540.444
$ws.Range("H132").Value = 2692.5745
$ws.Range("I132").Value = 1917.0312
$ws.Range("J132").Value = 4347.067
$ws.Range("K132").Value = 5751.0936
$ws.Range("L132").Value = 13041.201
$ws.Range("M132").Value = -3221.0936
$ws.Range("N132").Value = -18101.201
$ws.Range("H137").Value = 51151.6
$ws.Range("J137").Value = 51151.6
$ws.Range("L137").Value = 51151.6
$ws.Range("N137").Value = -61351.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 40000
$ws.Range("J98").Value = 40000
$ws.Range("L98").Value = 40000
$ws.Range("N98").Value = -45990
$ws.Range("H137").Value = 32945
$ws.Range("J137").Value = 32945
$ws.Range("L137").Value = 32945
$ws.Range("N137").Value = -43145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3289.1052
$ws.Range("I31").Value = 1347.6842
$ws.Range("J31").Value = 5230.5264
$ws.Range("K31").Value = 1347.6842
$ws.Range("L31").Value = 5230.5264
$ws.Range("M31").Value = -1052.6842
$ws.Range("N31").Value = -5820.5264
$ws.Range("H34").Value = 3289.1052
$ws.Range("I34").Value = 1347.6842
$ws.Range("J34").Value = 5230.5264
$ws.Range("K34").Value = 1347.6842
$ws.Range("L34").Value = 5230.5264
$ws.Range("M34").Value = -1145.6842
$ws.Range("N34").Value = -5634.5264
$ws.Range("H58").Value = 1938.8677
$ws.Range("I58").Value = 1689.0339
$ws.Range("K58").Value = 1689.0339
$ws.Range("M58").Value = -1486.0339
$ws.Range("H99").Value = 6253156
$ws.Range("I99").Value = 9525829
$ws.Range("J99").Value = 5326
$ws.Range("K99").Value = 9525829
$ws.Range("L99").Value = 5326
$ws.Range("M99").Value = -9524331
$ws.Range("N99").Value = -8322
$ws.Range("H107").Value = 668.8261
$ws.Range("I107").Value = 542.0476
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 542.0476
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1377.9524
$ws.Range("N107").Value = -5840
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H126").Value = 6253156
$ws.Range("I126").Value = 9525829
$ws.Range("J126").Value = 5326
$ws.Range("K126").Value = 28577487
$ws.Range("L126").Value = 15978
$ws.Range("M126").Value = -28575017
$ws.Range("N126").Value = -20918
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H130").Value = 41780
$ws.Range("J130").Value = 41780
$ws.Range("L130").Value = 41780
$ws.Range("N130").Value = -51820
$ws.Range("H132").Value = 3388
$ws.Range("I132").Value = 2580.6667
$ws.Range("J132").Value = 4670.2354
$ws.Range("K132").Value = 7742.000100000001
$ws.Range("L132").Value = 14010.7062
$ws.Range("M132").Value = -5212.000100000001
$ws.Range("N132").Value = -19070.7062
$ws.Range("H136").Value = 1938.8677
$ws.Range("I136").Value = 1689.0339
$ws.Range("K136").Value = 5067.101699999999
$ws.Range("M136").Value = -2517.101699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3482
$ws.Range("I87").Value = 3482
$ws.Range("K87").Value = 10446
$ws.Range("M87").Value = -9198
$ws.Range("H90").Value = 3482
$ws.Range("I90").Value = 3482
$ws.Range("K90").Value = 31338
$ws.Range("M90").Value = -25098
$ws.Range("H103").Value = 2249.75
$ws.Range("I103").Value = 1799.6
$ws.Range("J103").Value = 3000
$ws.Range("K103").Value = 5398.799999999999
$ws.Range("L103").Value = 9000
$ws.Range("M103").Value = -4519.799999999999
$ws.Range("N103").Value = -10758
$ws.Range("H131").Value = 10874857
$ws.Range("J131").Value = 785.4286
$ws.Range("L131").Value = 2356.2858
$ws.Range("N131").Value = -12436.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 20000
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20592
$ws.Range("H70").Value = 6362.683
$ws.Range("I70").Value = 5673.3335
$ws.Range("J70").Value = 9206.25
$ws.Range("K70").Value = 5673.3335
$ws.Range("L70").Value = 9206.25
$ws.Range("M70").Value = -5403.3335
$ws.Range("N70").Value = -9746.25
$ws.Range("H73").Value = 6362.683
$ws.Range("I73").Value = 5673.3335
$ws.Range("J73").Value = 9206.25
$ws.Range("K73").Value = 5673.3335
$ws.Range("L73").Value = 9206.25
$ws.Range("M73").Value = -4737.3335
$ws.Range("N73").Value = -11078.25
$ws.Range("H126").Value = 4161.7285
$ws.Range("I126").Value = 2822.5
$ws.Range("J126").Value = 5468.2925
$ws.Range("K126").Value = 8467.5
$ws.Range("L126").Value = 16404.8775
$ws.Range("M126").Value = -5997.5
$ws.Range("N126").Value = -21344.8775
$ws.Range("H137").Value = 72026.44500000001
$ws.Range("J137").Value = 72026.44500000001
$ws.Range("L137").Value = 72026.44500000001
$ws.Range("N137").Value = -82226.44500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1979.1666
$ws.Range("J46").Value = 2612.5
$ws.Range("L46").Value = 2612.5
$ws.Range("N46").Value = -2988.5
$ws.Range("H55").Value = 275.33334
$ws.Range("I55").Value = 202.71428
$ws.Range("K55").Value = 202.71428
$ws.Range("M55").Value = -29.71428
$ws.Range("H122").Value = 5581.1816
$ws.Range("I122").Value = 3299.125
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 9897.375
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -7447.375
$ws.Range("N122").Value = -39900.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 38000
$ws.Range("J106").Value = 38000
$ws.Range("L106").Value = 38000
$ws.Range("N106").Value = -40524
$ws.Range("H109").Value = 27377
$ws.Range("J109").Value = 27377
$ws.Range("L109").Value = 27377
$ws.Range("N109").Value = -30151
$ws.Range("H113").Value = 156.29411
$ws.Range("I113").Value = 150.93333
$ws.Range("J113").Value = 196.5
$ws.Range("K113").Value = 452.79999
$ws.Range("L113").Value = 589.5
$ws.Range("M113").Value = 1717.20001
$ws.Range("N113").Value = -4929.5

# These two cells were removed entirely (no LeveProfitHQ for an HQ-less recipe)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N123").ClearContents()
$ws.Range("N127").ClearContents()

Write-Output "Applied 203 cell updates and 2 cell removals"
